# Auto-generated edit script: updates Leve profit/price columns (H-N)
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled
# market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1198651.5
$ws.Range("I33").Value = 1278501.5
$ws.Range("K33").Value = 1278501.5
$ws.Range("M33").Value = -1278272.5

# Row 75
$ws.Range("H75").Value = 38400
$ws.Range("J75").Value = 38400
$ws.Range("L75").Value = 38400
$ws.Range("N75").Value = -40272

# Row 76
$ws.Range("H76").Value = 34000.25
$ws.Range("I76").Value = 65000.5
$ws.Range("K76").Value = 65000.5
$ws.Range("M76").Value = -64685.5

# Row 78
$ws.Range("H78").Value = 38400
$ws.Range("J78").Value = 38400
$ws.Range("L78").Value = 115200
$ws.Range("N78").Value = -124560

# Row 79
$ws.Range("H79").Value = 34000.25
$ws.Range("I79").Value = 65000.5
$ws.Range("K79").Value = 65000.5
$ws.Range("M79").Value = -63908.5

# Row 98
$ws.Range("H98").Value = 1736.1875
$ws.Range("I98").Value = 1198.5
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 1198.5
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = 299.5
$ws.Range("N98").Value = -8496

# Row 122
$ws.Range("H122").Value = 1736.1875
$ws.Range("I122").Value = 1198.5
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 3595.5
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -1145.5
$ws.Range("N122").Value = -21400

# Row 135
$ws.Range("H135").Value = 36934.215
$ws.Range("I135").Value = 46504.547
$ws.Range("K135").Value = 418540.923
$ws.Range("M135").Value = -416005.923

# Row 137
$ws.Range("H137").Value = 3572627.5
$ws.Range("I137").Value = 1614085.8
$ws.Range("J137").Value = 9092155
$ws.Range("K137").Value = 4842257.4
$ws.Range("L137").Value = 27276465
$ws.Range("M137").Value = -4839707.4
$ws.Range("N137").Value = -27281565

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2354.3547
$ws.Range("I2").Value = 1550.1786
$ws.Range("J2").Value = 9860
$ws.Range("K2").Value = 1550.1786
$ws.Range("L2").Value = 9860
$ws.Range("M2").Value = -1437.1786
$ws.Range("N2").Value = -10086

# Row 63
$ws.Range("H63").Value = 5500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 5500
$ws.Range("N63").Value = -6872
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 5500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 27500
$ws.Range("N66").Value = -34364
$ws.Range("M66").ClearContents()

# Row 88
$ws.Range("H88").Value = 40995.8
$ws.Range("I88").Value = 1623.6666
$ws.Range("J88").Value = 100054
$ws.Range("K88").Value = 1623.6666
$ws.Range("L88").Value = 100054
$ws.Range("M88").Value = -1217.6666
$ws.Range("N88").Value = -100866

# Row 91
$ws.Range("H91").Value = 40995.8
$ws.Range("I91").Value = 1623.6666
$ws.Range("J91").Value = 100054
$ws.Range("K91").Value = 1623.6666
$ws.Range("L91").Value = 100054
$ws.Range("M91").Value = -219.6666
$ws.Range("N91").Value = -102862

# Row 97
$ws.Range("H97").Value = 3070
$ws.Range("I97").Value = 1648.1818
$ws.Range("J97").Value = 5025
$ws.Range("K97").Value = 1648.1818
$ws.Range("L97").Value = 5025
$ws.Range("M97").Value = -1152.1818
$ws.Range("N97").Value = -6017

# Row 116
$ws.Range("H116").Value = 2354.3547
$ws.Range("I116").Value = 1550.1786
$ws.Range("J116").Value = 9860
$ws.Range("K116").Value = 1550.1786
$ws.Range("L116").Value = 9860
$ws.Range("M116").Value = 743.8214
$ws.Range("N116").Value = -14448

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2354.3547
$ws.Range("I3").Value = 1550.1786
$ws.Range("J3").Value = 9860
$ws.Range("K3").Value = 1550.1786
$ws.Range("L3").Value = 9860
$ws.Range("M3").Value = -1436.1786
$ws.Range("N3").Value = -10088

# Row 86
$ws.Range("H86").Value = 1891.5
$ws.Range("I86").Value = 1937.862
$ws.Range("J86").Value = 1699.4286
$ws.Range("K86").Value = 1937.862
$ws.Range("L86").Value = 1699.4286
$ws.Range("M86").Value = -814.8620000000001
$ws.Range("N86").Value = -3945.4286

# Row 89
$ws.Range("H89").Value = 1891.5
$ws.Range("I89").Value = 1937.862
$ws.Range("J89").Value = 1699.4286
$ws.Range("K89").Value = 9689.310000000001
$ws.Range("L89").Value = 8497.143
$ws.Range("M89").Value = -4073.310000000001
$ws.Range("N89").Value = -19729.143

# Row 94
$ws.Range("H94").Value = 882.2632
$ws.Range("I94").Value = 885.25
$ws.Range("J94").Value = 866.3333
$ws.Range("K94").Value = 885.25
$ws.Range("L94").Value = 866.3333
$ws.Range("M94").Value = -434.25
$ws.Range("N94").Value = -1768.3333

# Row 105
$ws.Range("H105").Value = 2746.0833
$ws.Range("I105").Value = 1695.3
$ws.Range("K105").Value = 1695.3
$ws.Range("M105").Value = 51.70000000000005

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3740
$ws.Range("I62").Value = 3842.8572
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3842.8572
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -3218.8572
$ws.Range("N62").Value = -4748

# Row 65
$ws.Range("H65").Value = 3740
$ws.Range("I65").Value = 3842.8572
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 19214.286
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -16094.286
$ws.Range("N65").Value = -23740

# Row 87
$ws.Range("H87").Value = 21598
$ws.Range("J87").Value = 21598
$ws.Range("L87").Value = 21598
$ws.Range("N87").Value = -23970

# Row 90
$ws.Range("H90").Value = 21598
$ws.Range("J90").Value = 21598
$ws.Range("L90").Value = 64794
$ws.Range("N90").Value = -76650

# Row 132
$ws.Range("H132").Value = 1830.0588
$ws.Range("I132").Value = 1688.8572
$ws.Range("K132").Value = 5066.571599999999
$ws.Range("M132").Value = -2536.571599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 20805.922
$ws.Range("I70").Value = 26319.938
$ws.Range("J70").Value = 4263.875
$ws.Range("K70").Value = 26319.938
$ws.Range("L70").Value = 4263.875
$ws.Range("M70").Value = -26049.938
$ws.Range("N70").Value = -4803.875

# Row 73
$ws.Range("H73").Value = 20805.922
$ws.Range("I73").Value = 26319.938
$ws.Range("J73").Value = 4263.875
$ws.Range("K73").Value = 26319.938
$ws.Range("L73").Value = 4263.875
$ws.Range("M73").Value = -25383.938
$ws.Range("N73").Value = -6135.875

# Row 80
$ws.Range("H80").Value = 1325486.5
$ws.Range("I80").Value = 2676.2307
$ws.Range("J80").Value = 2553810.5
$ws.Range("K80").Value = 2676.2307
$ws.Range("L80").Value = 2553810.5
$ws.Range("M80").Value = -1678.2307
$ws.Range("N80").Value = -2555806.5

# Row 83
$ws.Range("H83").Value = 1325486.5
$ws.Range("I83").Value = 2676.2307
$ws.Range("J83").Value = 2553810.5
$ws.Range("K83").Value = 13381.1535
$ws.Range("L83").Value = 12769052.5
$ws.Range("M83").Value = -8389.1535
$ws.Range("N83").Value = -12779036.5

# Row 123
$ws.Range("H123").Value = 34290
$ws.Range("J123").Value = 34290
$ws.Range("L123").Value = 34290
$ws.Range("N123").Value = -39190

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2154.3635
$ws.Range("I40").Value = 1671.4286
$ws.Range("J40").Value = 2999.5
$ws.Range("K40").Value = 1671.4286
$ws.Range("L40").Value = 2999.5
$ws.Range("M40").Value = -1535.4286
$ws.Range("N40").Value = -3271.5

# Row 82
$ws.Range("H82").Value = 1577
$ws.Range("I82").Value = 1384.2
$ws.Range("J82").Value = 1898.3334
$ws.Range("K82").Value = 1384.2
$ws.Range("L82").Value = 1898.3334
$ws.Range("M82").Value = -1023.2
$ws.Range("N82").Value = -2620.3334

# Row 85
$ws.Range("H85").Value = 1577
$ws.Range("I85").Value = 1384.2
$ws.Range("J85").Value = 1898.3334
$ws.Range("K85").Value = 1384.2
$ws.Range("L85").Value = 1898.3334
$ws.Range("M85").Value = -136.2
$ws.Range("N85").Value = -4394.3334

# Row 122
$ws.Range("H122").Value = 2675.3914
$ws.Range("I122").Value = 2076.353
$ws.Range("J122").Value = 3026.5518
$ws.Range("K122").Value = 6229.059
$ws.Range("L122").Value = 9079.6554
$ws.Range("M122").Value = -3779.059
$ws.Range("N122").Value = -13979.6554

# Row 132
$ws.Range("H132").Value = 1737.6666
$ws.Range("I132").Value = 1448.6888
$ws.Range("K132").Value = 4346.0664
$ws.Range("M132").Value = -1816.0664

# Row 136
$ws.Range("H136").Value = 1324.1642
$ws.Range("I136").Value = 1112.3148
$ws.Range("K136").Value = 3336.9444
$ws.Range("M136").Value = -786.9444000000003

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 10461.32
$ws.Range("I122").Value = 12027.579
$ws.Range("J122").Value = 5501.5
$ws.Range("K122").Value = 36082.737
$ws.Range("L122").Value = 16504.5
$ws.Range("M122").Value = -33632.737
$ws.Range("N122").Value = -21404.5

# Row 136
$ws.Range("H136").Value = 1202.6515
$ws.Range("I136").Value = 1205.7693
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3617.3079
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -1067.3079
$ws.Range("N136").Value = -8100
